$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "RATIOS" and "Rama descripción" columns (AU, BH) are being reclassified
# from iaest-measure to iaest-dimension, changing their identifier (row 3),
# type (row 4: medida -> dim), datatype (row 5: xsd:string -> skos:Concept),
# and adding a new metadata row (row 6) that references their mapping files.

# Row 3: identifier strings
$ws.Range("AU3").Value = "iaest-dimension:ratios"
$ws.Range("BH3").Value = "iaest-dimension:rama-descripcion"

# Row 4: medida -> dim
$ws.Range("AU4").Value = "dim"
$ws.Range("BH4").Value = "dim"

# Row 5: xsd:string -> skos:Concept
$ws.Range("AU5").Value = "skos:Concept"
$ws.Range("BH5").Value = "skos:Concept"

# Row 6 (new): mapping file references.
# Set the values first, then copy the row-5 cell formatting (style s="1")
# onto the new cells so they match the rest of the sheet's formatting.
$ws.Range("AU6").Value = "mapping-ratios.xlsx"
$ws.Range("BH6").Value = "mapping-rama-descripcion.xlsx"

$ws.Range("AU5").Copy()
$ws.Range("AU6").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("BH5").Copy()
$ws.Range("BH6").PasteSpecial(-4122)  # xlPasteFormats
